# Refresh the crypto snapshot: updated Price (D) and Volume(1h) (E) columns
# for the coin rows, per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "24.638.37"; E = "  -1.17%  " }
    @{ Row = 3; D = "1.672.74"; E = "  -2.10%  " }
    @{ Row = 4; D = "1.002"; E = "  -0.18%  " }
    @{ Row = 5; D = "313.49"; E = "  -0.66%  " }
    @{ Row = 6; D = $null; E = "  -0.17%  " }
    @{ Row = 7; D = "0.3936"; E = "  -2.18%  " }
    @{ Row = 8; D = "0.3945"; E = "  -3.12%  " }
    @{ Row = 9; D = "1.002"; E = "  -0.06%  " }
    @{ Row = 10; D = "1.399"; E = "  -5.56%  " }
    @{ Row = 11; D = "50.89"; E = "  -5.51%  " }
    @{ Row = 12; D = "0.08641"; E = "  -2.09%  " }
    @{ Row = 13; D = "25.20"; E = "  -4.82%  " }
    @{ Row = 14; D = "7.324"; E = "  -2.36%  " }
    @{ Row = 15; D = "0.00001316"; E = "  -2.17%  " }
    @{ Row = 16; D = "7.696"; E = "  -4.25%  " }
    @{ Row = 17; D = "1.677.22"; E = "  -2.00%  " }
    @{ Row = 18; D = "93.98"; E = "  -1.40%  " }
    @{ Row = 19; D = "0.07015"; E = "  -2.25%  " }
    @{ Row = 20; D = "21.07"; E = "  +0.57%  " }
    @{ Row = 21; D = "7.077"; E = "  -2.65%  " }
    @{ Row = 22; D = $null; E = "  -0.31%  " }
    @{ Row = 23; D = "13.93"; E = "  -3.98%  " }
    @{ Row = 24; D = "24.633.23"; E = "  -1.14%  " }
    @{ Row = 25; D = "2.348"; E = "  +0.26%  " }
    @{ Row = 26; D = "2.785"; E = "  -3.85%  " }
    @{ Row = 27; D = "23.00"; E = "  -0.68%  " }
    @{ Row = 28; D = "5.834"; E = "  -8.60%  " }
    @{ Row = 29; D = "160.07"; E = "  -2.06%  " }
    @{ Row = 30; D = $null; E = "  +1.55%  " }
    @{ Row = 31; D = "8.338"; E = "  +1.56%  " }
    @{ Row = 32; D = "2.485"; E = "  +9.54%  " }
    @{ Row = 33; D = "1.863.67"; E = "  -0.53%  " }
    @{ Row = 34; D = "0.03082"; E = "  -3.78%  " }
    @{ Row = 35; D = "0.08253"; E = "  -5.65%  " }
    @{ Row = 36; D = "6.957"; E = "  -5.26%  " }
    @{ Row = 37; D = "0.2813"; E = "  -2.23%  " }
    @{ Row = 38; D = "0.9915"; E = "  -3.81%  " }
    @{ Row = 39; D = "0.09593"; E = "  +1.35%  " }
    @{ Row = 40; D = "1.516"; E = "  +2.62%  " }
    @{ Row = 41; D = "10.30"; E = "  -5.40%  " }
    @{ Row = 42; D = "0.7903"; E = "  -6.82%  " }
    @{ Row = 43; D = "13.48"; E = "  -4.82%  " }
    @{ Row = 44; D = "16.62"; E = "  -5.74%  " }
    @{ Row = 45; D = "2.562"; E = $null }
    @{ Row = 46; D = "0.7094"; E = "  -5.03%  " }
    @{ Row = 47; D = "4.167"; E = "  -1.61%  " }
    @{ Row = 48; D = "0.08669"; E = "  +2.97%  " }
    @{ Row = 49; D = "1.001"; E = "  -0.17%  " }
    @{ Row = 50; D = "1.327"; E = "  -4.67%  " }
    @{ Row = 51; D = "137.90"; E = "  -2.48%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        # Force text so price strings like "25.20" or "1.001" keep their
        # exact digits/trailing zeros instead of being parsed as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
